# Weekly fruit/vegetable price update:
# a new observation is inserted as row 194 (pushing the existing rows
# 194..227 down to 195..228), matching the source data feed's ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 194, shifting rows 194:227 down to 195:228
# (carries formatting, e.g. the date style on column D, down with it).
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with the new price observation.
$ws.Cells.Item(194, 1).Value  = 5
$ws.Cells.Item(194, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(194, 3).Value  = "Maule"
$ws.Cells.Item(194, 4).Value  = 44984
$ws.Cells.Item(194, 5).Value  = 7
$ws.Cells.Item(194, 6).Value  = 100112031
$ws.Cells.Item(194, 7).Value  = "Poroto verde"
$ws.Cells.Item(194, 8).Value  = "Sin especificar"
$ws.Cells.Item(194, 9).Value  = "Primera"
$ws.Cells.Item(194, 10).Value = 150
$ws.Cells.Item(194, 11).Value = 25000
$ws.Cells.Item(194, 12).Value = 25000
$ws.Cells.Item(194, 13).Value = 25000
$ws.Cells.Item(194, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(194, 15).Value = "Región del Maule"
$ws.Cells.Item(194, 16).Value = 1000
$ws.Cells.Item(194, 17).Value = 25
$ws.Cells.Item(194, 18).Value = "Hortaliza"
